$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "320018594180",
    "320018594190",
    "320018594227",
    "320018594249",
    "320018594282",
    "320018594308",
    "320018594330",
    "320018594352",
    "320018594385",
    "320018594400",
    "320018594444",
    "320018594466",
    "320018594499",
    "320018594514",
    "320018594547",
    "320018594569",
    "320018594606",
    "320018594628",
    "320018594650",
    "320018594672",
    "320018594709"
)

$mirrorRows = @(5, 6, 7, 13, 14, 15, 16, 17)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $val = $values[$i]

    $cC = $ws.Cells.Item($row, 3)
    $cC.Value = "'" + $val
    $cC.Style = "Normal"

    if ($mirrorRows -contains $row) {
        $cD = $ws.Cells.Item($row, 4)
        $cD.Value = "'" + $val
        $cD.Style = "Normal"
    }
}
